# Weekly fruit/vegetable price update: a new daily record needs to be
# inserted into the "Coliflor" price history table, shifting the existing
# rows down by one (the sheet keeps growing with one new logged price per
# update). The new record is inserted at row 212; everything at or below
# that row moves down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 212 - this pushes the previous rows
# 212..334 down to 213..335 and extends the used range to row 335.
$ws.Rows(212).Insert()

# Populate the newly inserted row with the new weekly observation.
$ws.Range("A212").Value2 = 4
$ws.Range("B212").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C212").Value2 = "Los Lagos"
$ws.Range("D212").Value2 = 44719
$ws.Range("E212").Value2 = 10
$ws.Range("F212").Value2 = 100112008
$ws.Range("G212").Value2 = "Coliflor"
$ws.Range("H212").Value2 = "Sin especificar"
$ws.Range("I212").Value2 = "Primera"
$ws.Range("J212").Value2 = 1000
$ws.Range("K212").Value2 = 1600
$ws.Range("L212").Value2 = 1700
$ws.Range("M212").Value2 = 1650
$ws.Range("N212").Value2 = "$/unidad"
$ws.Range("O212").Value2 = "Región Metropolitana"
$ws.Range("P212").Value2 = 1650
$ws.Range("Q212").Value2 = 1
$ws.Range("R212").Value2 = "Hortaliza"
